# Apply cell updates per diff (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.554.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.818.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "665.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.816.33"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.51%  "
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.462.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.827.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.545.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "476.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("E23").Value = "  +1.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.971.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("E31").Value = "  +7.63%  "
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.180"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.28%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.13%  "
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.776.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("E40").Value = "  +1.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.973"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E45").Value = "  +9.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "158.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("E50").Value = "  +4.31%  "
$ws.Range("E51").Value = "  +1.29%  "
